$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "24/10/2025"
$ws.Range("B14").Value = "Aarhus"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = "Nordsjaelland"
$ws.Range("F14").Value = "W"
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1.64
$ws.Range("L14").Value = 0.63
$ws.Range("M14").Value = 24
$ws.Range("N14").Value = 6
$ws.Range("O14").Value = 8
$ws.Range("P14").Value = 3
